$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column P (year 2022) mirroring the formatting of column O.

# Row 3: empty bottom-border cell, same style as O3.
$ws.Range("O3").Copy()
$ws.Range("P3").PasteSpecial(-4122)
$ws.Range("P3").ClearContents()

# Row 4: header year value 2022, same style as O4.
$ws.Range("O4").Copy()
$ws.Range("P4").PasteSpecial(-4122)
$ws.Range("P4").Value = 2022

# Row 5: empty cell, same style as O5.
$ws.Range("O5").Copy()
$ws.Range("P5").PasteSpecial(-4122)
$ws.Range("P5").ClearContents()

# Row 6: data value 1373, same style as O6.
$ws.Range("O6").Copy()
$ws.Range("P6").PasteSpecial(-4122)
$ws.Range("P6").Value = 1373

# Row 7: text value "-", same style as O7 plus right alignment.
$ws.Range("O7").Copy()
$ws.Range("P7").PasteSpecial(-4122)
$ws.Range("P7").Value = "-"
$ws.Range("P7").HorizontalAlignment = -4152

# Row 8: data value 117, same style as O8.
$ws.Range("O8").Copy()
$ws.Range("P8").PasteSpecial(-4122)
$ws.Range("P8").Value = 117

# Row 9: data value 154, same style as O9.
$ws.Range("O9").Copy()
$ws.Range("P9").PasteSpecial(-4122)
$ws.Range("P9").Value = 154

# Row 10: data value 885, same style as O10.
$ws.Range("O10").Copy()
$ws.Range("P10").PasteSpecial(-4122)
$ws.Range("P10").Value = 885

# Move the active selection to P7, matching the saved view state.
$ws.Range("P7").Select()
